$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "release" values (from home) for columns B (Atkin), E (Eratosthenes), H (Sundaram)
$values = @{
    4  = @(355156, 332321, 363086)
    5  = @(358755, 333326, 375921)
    6  = @(355651, 330839, 373563)
    7  = @(356771, 333033, 384730)
    8  = @(356573, 334187, 382976)
    9  = @(356471, 331746, 379426)
    10 = @(357041, 331386, 399030)
    11 = @(356491, 331100, 353635)
    12 = @(356949, 332007, 361482)
    13 = @(359119, 329878, 380121)
}

foreach ($row in $values.Keys) {
    $trio = $values[$row]

    $ws.Range("B$row").Value = $trio[0]
    $ws.Range("E$row").Value = $trio[1]

    $ws.Range("H$row").Value = $trio[2]
    # The H column loses its centered "s=2" style in the target workbook,
    # reverting back to the default "Normal" style.
    $ws.Range("H$row").Style = "Normal"
}

# Update the selection shown in the saved sheet view (was R4:R13, now B4:B13)
$ws.Range("B4:B13").Select()
